$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample number text from E7420 to E7420L (shared string used in column G)
$ws.Range("G2:G37").Value = "E7420L"

# Replace the FALSE() formula cells in column H with plain boolean FALSE values
$ws.Range("H2:H37").Value = $false
